$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Бухгалтерия" (accounting) use-case rows that described
# "Формирование персональных скидок" (row 16) and "Формирование зарплат"
# (row 17). Deleting the entire rows shifts everything below them up by two
# rows (old rows 18-19 become new rows 16-17), and the sheet's used range
# shrinks from A1:E19 to A1:E17.
$ws.Range("A16:E17").EntireRow.Delete()

# Mirror the selection Excel leaves on the rows that used to hold the
# deleted records.
$ws.Range("A16:XFD17").Select()
